# Update forecast values in column C to reflect refreshed model output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = 41.90542479266645
$ws.Range("C3").Value  = 38.49402873427938
$ws.Range("C4").Value  = 6.425165193507417
$ws.Range("C5").Value  = 47.09297571453026
$ws.Range("C6").Value  = 40.73434009058274
$ws.Range("C7").Value  = 47.95405382554818
$ws.Range("C8").Value  = 44.18388664128612

$ws.Range("C10").Value = 36.05412462389766
$ws.Range("C11").Value = 22.6026861500278
$ws.Range("C12").Value = 1.500250415078458
$ws.Range("C13").Value = 0
$ws.Range("C15").Value = 4.65795748237634

$ws.Range("C93").Value  = 0.08146223558543657
$ws.Range("C94").Value  = 0.01615361644590407
$ws.Range("C95").Value  = 2.097005524008066
$ws.Range("C96").Value  = 0.002014158141016677
$ws.Range("C97").Value  = 1.144186650192057
$ws.Range("C98").Value  = 1.259004036088275
$ws.Range("C99").Value  = 0
$ws.Range("C100").Value = 7.511973968952613
$ws.Range("C101").Value = 0.92962479075577
$ws.Range("C102").Value = 2.349444947102698
$ws.Range("C103").Value = 0.06719329937135948
$ws.Range("C104").Value = 0
$ws.Range("C105").Value = 0
$ws.Range("C106").Value = 4.827966511230345

$ws.Range("C107").Value = 1.509224491614277
$ws.Range("C108").Value = 12.95482561237078
$ws.Range("C110").Value = 0.71157218845641
$ws.Range("C112").Value = 0.1952368699170393
$ws.Range("C113").Value = 3.248394638311456

$ws.Range("C121").Value = 9.521895170679079
$ws.Range("C122").Value = 4.056154768915386
$ws.Range("C123").Value = 0
$ws.Range("C124").Value = 0
$ws.Range("C125").Value = 179.8432663935019
$ws.Range("C126").Value = 3.034031357935503
$ws.Range("C127").Value = 63.31538421522037
